# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - Salario Basico actualizado de 1.000.000 a 3.000.000 para los periodos
#   existentes (2506, 2507) y el nuevo periodo agregado (2508).
# - Se agrega un nuevo periodo de mora (2508) para el mismo trabajador,
#   lo que desplaza hacia abajo el bloque de firma.
# - Totales (Valor Mora / Cant. Periodos) actualizados de acuerdo al nuevo
#   periodo agregado.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserta una fila nueva para el tercer periodo (2508), lo que empuja hacia
# abajo el bloque de firma (filas 22-23 -> 23-24).
$ws.Rows.Item(18).Insert()

# La fila 17 (antigua "ultima fila" con borde inferior solido) se convierte
# en la fila intermedia: copia el formato de la fila 18 recien insertada
# para que tenga el borde solido de "ultima fila" de la tabla.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

# La fila 17 pasa a ser una fila intermedia: copia el formato de la fila 16.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# --- Totales ---
$ws.Range("E11").Value2 = 360000   # VALOR MORA
$ws.Range("F13").Value2 = 3        # Cant. Periodos

# --- Fila 16: periodo 2506 ---
$ws.Range("E16").Value2 = "2506"
$ws.Range("G16").Value2 = 3000000

# --- Fila 17: periodo 2507 ---
$ws.Range("E17").Value2 = "2507"
$ws.Range("G17").Value2 = 3000000

# --- Fila 18 (nueva): periodo 2508 ---
$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "22800731"
$ws.Range("D18").Value2 = "CARMEN JULIA GUERRERO CALLE"
$ws.Range("E18").Value2 = "2508"
$ws.Range("F18").Value2 = 120000
$ws.Range("G18").Value2 = 3000000

$wb.Save()
